$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.842.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.786.75'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.16%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '443.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +12.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.620'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.92%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +2.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.150'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000308'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.86'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.30'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.376.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.99%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.853.52'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.79%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.137'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('E19').Value = '  +6.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '66.864.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '411.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.49'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '36.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.39'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +29.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.67'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '730.43'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.88'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +13.69%  '
$ws.Range('E32').Value = '  +11.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.36'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +15.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.158'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '56.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  +24.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0474'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +30.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.88'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.338'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +18.56%  '
$ws.Range('E43').Value = '  +4.79%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.40%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₃0668'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -11.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.08'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.64'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '143.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.82'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.08%  '
